$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.868.76"
$ws.Range("E2").Value = "  +4.54%  "

# Row 3
$ws.Range("D3").Value = "2.378.01"
$ws.Range("E3").Value = "  +3.20%  "

# Row 4
$ws.Range("E4").Value = "  -1.09%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "317.35"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "111.08"
$ws.Range("E6").Value = "  +6.29%  "

# Row 7
$ws.Range("E7").Value = "  +1.17%  "

# Row 8
$ws.Range("E8").Value = "  -0.55%  "

# Row 9
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  +3.94%  "

# Row 10
$ws.Range("D10").Value = "42.14"
$ws.Range("E10").Value = "  +6.78%  "

# Row 11
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  +2.34%  "

# Row 12
$ws.Range("D12").Value = "8.67"
$ws.Range("E12").Value = "  +4.15%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +1.57%  "

# Row 15
$ws.Range("D15").Value = "15.72"
$ws.Range("E15").Value = "  +3.38%  "

# Row 16
$ws.Range("D16").Value = "2.739.75"
$ws.Range("E16").Value = "  +3.16%  "

# Row 17
$ws.Range("D17").Value = "2.370.76"
$ws.Range("E17").Value = "  +2.82%  "

# Row 18
$ws.Range("D18").Value = "44.817.48"
$ws.Range("E18").Value = "  +4.64%  "

# Row 19
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +3.86%  "

# Row 20
$ws.Range("E20").Value = "  +2.39%  "

# Row 21
$ws.Range("D21").Value = "13.23"
$ws.Range("E21").Value = "  -3.09%  "

# Row 22
$ws.Range("D22").Value = "75.04"
$ws.Range("E22").Value = "  +2.30%  "

# Row 23
$ws.Range("E23").Value = "  +1.24%  "

# Row 24
$ws.Range("D24").Value = "267.54"
$ws.Range("E24").Value = "  +1.40%  "

# Row 25
$ws.Range("E25").Value = "  +6.04%  "

# Row 26
$ws.Range("E26").Value = "  -0.74%  "

# Row 27
$ws.Range("D27").Value = "7.76"
$ws.Range("E27").Value = "  +12.14%  "

# Row 28
$ws.Range("D28").Value = "11.22"
$ws.Range("E28").Value = "  +4.29%  "

# Row 29
$ws.Range("E29").Value = "  -0.25%  "

# Row 30
$ws.Range("D30").Value = "39.29"
$ws.Range("E30").Value = "  +8.12%  "

# Row 31
$ws.Range("D31").Value = "22.71"
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("D32").Value = "168.38"
$ws.Range("E32").Value = "  +1.48%  "

# Row 33
$ws.Range("D33").Value = "0.0915"
$ws.Range("E33").Value = "  +5.36%  "

# Row 34
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").Value = "  +12.68%  "

# Row 35
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +3.69%  "

# Row 37
$ws.Range("D37").Value = "4.79"
$ws.Range("E37").Value = "  +5.68%  "

# Row 38
$ws.Range("D38").Value = "0.0365"
$ws.Range("E38").Value = "  +4.02%  "

# Row 39
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  +9.10%  "

# Row 40
$ws.Range("D40").Value = "3.88"
$ws.Range("E40").Value = "  +3.11%  "

# Row 41
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  +9.06%  "

# Row 42
$ws.Range("D42").Value = "104.24"
$ws.Range("E42").Value = "  +5.43%  "

# Row 43
$ws.Range("D43").Value = "13.78"
$ws.Range("E43").Value = "  +13.24%  "

# Row 44
$ws.Range("E44").Value = "  +5.49%  "

# Row 45
$ws.Range("D45").Value = "71.17"
$ws.Range("E45").Value = "  +1.81%  "

# Row 47
$ws.Range("D47").Value = "119.39"
$ws.Range("E47").Value = "  +6.73%  "

# Row 48
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "5.58"
$ws.Range("E48").Value = "  +6.87%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "79.42"
$ws.Range("E49").Value = "  -1.72%  "

# Row 50
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.221"
$ws.Range("E50").Value = "  +16.44%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "9.07"
$ws.Range("E51").Value = "  +4.13%  "
